$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    69558.77063519611,
    65789.22778384826,
    62805.47619618164,
    60663.98826309984,
    60623.33912578635,
    63261.62112043737,
    63872.87106878438,
    70311.25048947106,
    87888.86178402032,
    95639.15601608463,
    100425.8142991323,
    102357.4719302255,
    101962.3151288686,
    105820.1113070404,
    105624.6049056241,
    104640.2348783341,
    99698.72768430505,
    94801.62774983187,
    92087.17084853728,
    88084.3101515717,
    85698.29265509204,
    83514.34032068169,
    80108.54038067423,
    75633.44261675818
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 2).Value = $val
    $row++
}
